$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worker/period table (rows 16-31, columns B:G) is rebuilt: the four
# workers now repeat across ascending periods 2505, 2506, 2507, 2508
# (previously 2507, 2506, 2505, 2504), and the "LICED MARGARITA TAPIA
# TORRES" rows now carry the same 1423500 Valor Mora as everyone else
# (was 1900000).

$workers = @(
    @("CC", "73140520", "OSWALDO DE JESUS CASTILLA TARRA"),
    @("CC", "1102868229", "FREDY DE JESUS MENDOZA PEREZ"),
    @("CC", "1052079546", "MARIA JOSE GONZALEZ ANGULO"),
    @("CC", "45649374", "LICED MARGARITA TAPIA TORRES")
)

$periods = @("2505", "2506", "2507", "2508")

$row = 16
foreach ($period in $periods) {
    foreach ($worker in $workers) {
        $ws.Cells.Item($row, 2).Value = $worker[0]
        $ws.Cells.Item($row, 3).Value = $worker[1]
        $ws.Cells.Item($row, 4).Value = $worker[2]
        $ws.Cells.Item($row, 5).Value = $period
        $ws.Cells.Item($row, 6).Value = 56940
        $ws.Cells.Item($row, 7).Value = 1423500
        $row = $row + 1
    }
}
